# Automatische test-sync: 2025-06-23 18:31:50
# Adds the new "Retour aanmelden" log entry (row 14) to the Logs sheet,
# extends the conditional-formatting ranges to include it, and refreshes
# the Dashboard category summary (and its row order) to account for the
# new "Retour / Terugbetaling" entry.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- Logs!A14:G14 --------------------------------------------------------
$logs.Range("A14").Value = "Retour aanmelden"
$logs.Range("B14").Value = "mailmind.test@zohomail.eu"
$logs.Range("C14").Value = "Ik wil graag een artikel retourneren. Hoe werkt dat?"
$logs.Range("D14").Value = "Retour / Terugbetaling"
$logs.Range("E14").Value = "Beste klant,`nBedankt voor je bericht. Om een artikel te retourneren kun je het volgende stappenplan volgen:`n1. Neem contact met ons op via retour@bedrijfsnaam.nl met je bestelnummer en de reden van retourneren.`n2. Wij zullen je vervolgens instructies sturen over het retourproces en het retouradres.`n3. Zodra wij het geretourneerde artikel hebben ontvangen, zullen we de terugbetaling verwerken.`nMocht je nog verdere vragen hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Naam]  `nKlantenservice Team  `nBedrijfsnaam"
$logs.Range("F14").Value = "2025-06-23 18:30:51"
$logs.Range("G14").Value = "Ja"

# Row height auto-adjusts when a multi-line value (E14) is entered;
# re-run AutoFit so the row keeps the sheet's implicit default height,
# same as every other row in the log.
$logs.Rows.Item(14).AutoFit()

# --- extend conditional formatting to cover the new row ------------------
$catFormats = $logs.Range("D2:D13").FormatConditions
for ($i = 1; $i -le $catFormats.Count; $i++) {
    $catFormats.Item($i).ModifyAppliesToRange($logs.Range("D2:D14"))
}

$answeredFormats = $logs.Range("G2:G13").FormatConditions
for ($i = 1; $i -le $answeredFormats.Count; $i++) {
    $answeredFormats.Item($i).ModifyAppliesToRange($logs.Range("G2:G14"))
}

# --- Dashboard!A6:B7 summary reorder --------------------------------------
# "Retour / Terugbetaling" now has 2 occurrences and sorts above
# "Sollicitatie / Vacature" (still 1 occurrence).
$dash.Range("A6").Value = "Retour / Terugbetaling"
$dash.Range("B6").Value = 2
$dash.Range("A7").Value = "Sollicitatie / Vacature"
$dash.Range("B7").Value = 1
